# Apply updated TPM-derived values to Cxcl12-Itga4 LR-pair sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 146.620486
$ws.Range("H2").Value = 439.861458
$ws.Range("I2").Value = 0.3983053592962091
$ws.Range("J2").Value = 0.3983053592962091
$ws.Range("K2").Value = 3.0
$ws.Range("L2").Value = 1.0
$ws.Range("M2").Value = 0.4702473333333333
$ws.Range("N2").Value = 1.410742
$ws.Range("O2").Value = 0.00903492226842282
$ws.Range("P2").Value = 0.00903492226842282
$ws.Range("Q2").Value = 68.94789255353733
$ws.Range("R2").Value = 620.5310329818359
$ws.Range("S2").Value = 0.003598657960337471
$ws.Range("T2").Value = 0.003598657960337471
$ws.Range("G3").Value = 146.620486
$ws.Range("H3").Value = 439.861458
$ws.Range("I3").Value = 0.3983053592962091
$ws.Range("J3").Value = 0.3983053592962091
$ws.Range("N3").Value = 0.9584440000000001
$ws.Range("O3").Value = 0.006138235792679485
$ws.Range("P3").Value = 0.006138235792679485
$ws.Range("Q3").Value = 46.84250836126133
$ws.Range("R3").Value = 421.582575251352
$ws.Range("S3").Value = 0.002444892212848053
$ws.Range("T3").Value = 0.002444892212848053
$ws.Range("G4").Value = 146.620486
$ws.Range("H4").Value = 439.861458
$ws.Range("I4").Value = 0.3983053592962091
$ws.Range("J4").Value = 0.3983053592962091
$ws.Range("M4").Value = 1.047307
$ws.Range("N4").Value = 3.141921
$ws.Range("O4").Value = 0.02012204358311108
$ws.Range("P4").Value = 0.02012204358311108
$ws.Range("Q4").Value = 153.556661331202
$ws.Range("R4").Value = 1382.009951980818
$ws.Range("S4").Value = 0.008014717799145038
$ws.Range("T4").Value = 0.008014717799145038
$ws.Range("G5").Value = 146.620486
$ws.Range("H5").Value = 439.861458
$ws.Range("I5").Value = 0.3983053592962091
$ws.Range("J5").Value = 0.3983053592962091
$ws.Range("M5").Value = 50.21070966666667
$ws.Range("N5").Value = 150.632129
$ws.Range("O5").Value = 0.9647047983557866
$ws.Range("P5").Value = 0.9647047983557866
$ws.Range("Q5").Value = 7361.918653731565
$ws.Range("R5").Value = 66257.26788358409
$ws.Range("S5").Value = 0.3842470913238785
$ws.Range("T5").Value = 0.3842470913238785
$ws.Range("I6").Value = 0.534552907532962
$ws.Range("J6").Value = 0.5345529075329621
$ws.Range("K6").Value = 3.0
$ws.Range("L6").Value = 1.0
$ws.Range("M6").Value = 0.4702473333333333
$ws.Range("N6").Value = 1.410742
$ws.Range("O6").Value = 0.00903492226842282
$ws.Range("P6").Value = 0.00903492226842282
$ws.Range("Q6").Value = 92.53276555928689
$ws.Range("R6").Value = 832.794890033582
$ws.Range("S6").Value = 0.004829643967919723
$ws.Range("T6").Value = 0.004829643967919724
$ws.Range("I7").Value = 0.534552907532962
$ws.Range("J7").Value = 0.5345529075329621
$ws.Range("N7").Value = 0.9584440000000001
$ws.Range("O7").Value = 0.006138235792679485
$ws.Range("P7").Value = 0.006138235792679485
$ws.Range("Q7").Value = 62.86583510925823
$ws.Range("R7").Value = 565.792515983324
$ws.Range("S7").Value = 0.003281211790099715
$ws.Range("T7").Value = 0.003281211790099715
$ws.Range("I8").Value = 0.534552907532962
$ws.Range("J8").Value = 0.5345529075329621
$ws.Range("M8").Value = 1.047307
$ws.Range("N8").Value = 3.141921
$ws.Range("O8").Value = 0.02012204358311108
$ws.Range("P8").Value = 0.02012204358311108
$ws.Range("Q8").Value = 206.0834931538157
$ws.Range("R8").Value = 1854.751438384341
$ws.Range("S8").Value = 0.01075629690285701
$ws.Range("T8").Value = 0.01075629690285701
$ws.Range("I9").Value = 0.534552907532962
$ws.Range("J9").Value = 0.5345529075329621
$ws.Range("M9").Value = 50.21070966666667
$ws.Range("N9").Value = 150.632129
$ws.Range("O9").Value = 0.9647047983557866
$ws.Range("P9").Value = 0.9647047983557866
$ws.Range("Q9").Value = 9880.19600923008
$ws.Range("R9").Value = 88921.76408307072
$ws.Range("S9").Value = 0.5156857548720856
$ws.Range("T9").Value = 0.5156857548720857
$ws.Range("G10").Value = 24.174389
$ws.Range("H10").Value = 72.523167
$ws.Range("I10").Value = 0.0656715098899026
$ws.Range("J10").Value = 0.0656715098899026
$ws.Range("K10").Value = 3.0
$ws.Range("L10").Value = 1.0
$ws.Range("M10").Value = 0.4702473333333333
$ws.Range("N10").Value = 1.410742
$ws.Range("O10").Value = 0.00903492226842282
$ws.Range("P10").Value = 0.00903492226842282
$ws.Range("Q10").Value = 11.36794196221267
$ws.Range("R10").Value = 102.311477659914
$ws.Range("S10").Value = 0.0005933369871052304
$ws.Range("T10").Value = 0.0005933369871052304
$ws.Range("G11").Value = 24.174389
$ws.Range("H11").Value = 72.523167
$ws.Range("I11").Value = 0.0656715098899026
$ws.Range("J11").Value = 0.0656715098899026
$ws.Range("N11").Value = 0.9584440000000001
$ws.Range("O11").Value = 0.006138235792679485
$ws.Range("P11").Value = 0.006138235792679485
$ws.Range("Q11").Value = 7.723266030238667
$ws.Range("R11").Value = 69.509394272148
$ws.Range("S11").Value = 0.0004031072125655049
$ws.Range("T11").Value = 0.0004031072125655049
$ws.Range("G12").Value = 24.174389
$ws.Range("H12").Value = 72.523167
$ws.Range("I12").Value = 0.0656715098899026
$ws.Range("J12").Value = 0.0656715098899026
$ws.Range("M12").Value = 1.047307
$ws.Range("N12").Value = 3.141921
$ws.Range("O12").Value = 0.02012204358311108
$ws.Range("P12").Value = 0.02012204358311108
$ws.Range("Q12").Value = 25.318006820423
$ws.Range("R12").Value = 227.862061383807
$ws.Range("S12").Value = 0.001321444984173331
$ws.Range("T12").Value = 0.001321444984173331
$ws.Range("G13").Value = 24.174389
$ws.Range("H13").Value = 72.523167
$ws.Range("I13").Value = 0.0656715098899026
$ws.Range("J13").Value = 0.0656715098899026
$ws.Range("M13").Value = 50.21070966666667
$ws.Range("N13").Value = 150.632129
$ws.Range("O13").Value = 0.9647047983557866
$ws.Range("P13").Value = 0.9647047983557866
$ws.Range("Q13").Value = 1213.813227448061
$ws.Range("R13").Value = 10924.31904703254
$ws.Range("S13").Value = 0.06335362070605853
$ws.Range("T13").Value = 0.06335362070605853
$ws.Range("E14").Value = 3.0
$ws.Range("F14").Value = 1.0
$ws.Range("G14").Value = 0.541205
$ws.Range("H14").Value = 1.623615
$ws.Range("I14").Value = 0.001470223280926138
$ws.Range("J14").Value = 0.001470223280926138
$ws.Range("K14").Value = 3.0
$ws.Range("L14").Value = 1.0
$ws.Range("M14").Value = 0.4702473333333333
$ws.Range("N14").Value = 1.410742
$ws.Range("O14").Value = 0.00903492226842282
$ws.Range("P14").Value = 0.00903492226842282
$ws.Range("Q14").Value = 0.2545002080366667
$ws.Range("R14").Value = 2.29050187233
$ws.Range("S14").Value = 0.00001328335306039322
$ws.Range("T14").Value = 0.00001328335306039322
$ws.Range("E15").Value = 3.0
$ws.Range("F15").Value = 1.0
$ws.Range("G15").Value = 0.541205
$ws.Range("H15").Value = 1.623615
$ws.Range("I15").Value = 0.001470223280926138
$ws.Range("J15").Value = 0.001470223280926138
$ws.Range("N15").Value = 0.9584440000000001
$ws.Range("O15").Value = 0.006138235792679485
$ws.Range("P15").Value = 0.006138235792679485
$ws.Range("Q15").Value = 0.1729048950066667
$ws.Range("R15").Value = 1.55614405506
$ws.Range("S15").Value = 0.000009024577166211485
$ws.Range("T15").Value = 0.000009024577166211485
$ws.Range("E16").Value = 3.0
$ws.Range("F16").Value = 1.0
$ws.Range("G16").Value = 0.541205
$ws.Range("H16").Value = 1.623615
$ws.Range("I16").Value = 0.001470223280926138
$ws.Range("J16").Value = 0.001470223280926138
$ws.Range("M16").Value = 1.047307
$ws.Range("N16").Value = 3.141921
$ws.Range("O16").Value = 0.02012204358311108
$ws.Range("P16").Value = 0.02012204358311108
$ws.Range("Q16").Value = 0.5668077849350001
$ws.Range("R16").Value = 5.101270064415
$ws.Range("S16").Value = 0.00002958389693570032
$ws.Range("T16").Value = 0.00002958389693570032
$ws.Range("E17").Value = 3.0
$ws.Range("F17").Value = 1.0
$ws.Range("G17").Value = 0.541205
$ws.Range("H17").Value = 1.623615
$ws.Range("I17").Value = 0.001470223280926138
$ws.Range("J17").Value = 0.001470223280926138
$ws.Range("M17").Value = 50.21070966666667
$ws.Range("N17").Value = 150.632129
$ws.Range("O17").Value = 0.9647047983557866
$ws.Range("P17").Value = 0.9647047983557866
$ws.Range("Q17").Value = 27.17428712514834
$ws.Range("R17").Value = 244.568584126335
$ws.Range("S17").Value = 0.001418331453763833
$ws.Range("T17").Value = 0.001418331453763833
